$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")
$ws.Activate()

# --- Cell value updates (hours consumed per day) ---
# Row 52 - Codificacion (CU-01)
$ws.Range("Z52").Value = 3
$ws.Range("AC52").Value = 2
$ws.Range("AF52").Value = 5
$ws.Range("AI52").Value = 3

# Row 53 - Pruebas Unitarias (CU-01)
$ws.Range("AR53").Value = 1

# Row 55 - Pruebas Unitarias (CU-25)
$ws.Range("Q55").Value = 1.5

# Row 57 - Diagrama de robustez (CU-02)
$ws.Range("K57").Value = 1

# Row 58 - Diagrama de secuencia (CU-02)
$ws.Range("K58").Value = 1

# Row 60 - Diagrama de robustez (CU-03)
$ws.Range("N60").Value = 1

# Row 61 - Diagrama de secuencia (CU-03)
$ws.Range("N61").Value = 1

# --- Row height change on row 56 ---
$ws.Rows.Item(56).RowHeight = 48

# --- Update frozen pane / selection to reflect latest work location ---
$ws.Range("Y62").Select()
$excel.ActiveWindow.ScrollRow = 53
$excel.ActiveWindow.ScrollColumn = 43

$wb.Save()
